# Autogenerated on Thu Mar 26 2015 18:06:15 GMT+0000 (Coordinated Universal Time)
#
# Replaces the KOSIS-sourced citation block (A54:A61) on the "Summary" sheet
# with the Small Business Corporation citation block (A54:A62), splitting
# the old single-line source text into three lines (blank / long citation /
# URL), dropping the live hyperlink on the URL cell, and rewriting the final
# long citation string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the worksheet hyperlink that used to live on the URL line (old A56) -
# the rebuilt block no longer links the URL text.
foreach ($hl in @($ws.Hyperlinks)) {
    $hl.Delete()
}

# Clear out the old citation block (A54:A61) so it can be rebuilt with the
# new row layout (A54:A62).
$ws.Range("A54:A61").Clear()

$ws.Range("A54").Value = "Source:"
$ws.Range("A54").Style = "source"

$ws.Range("A55").Value = ""
$ws.Range("A55").Style = "source"

$ws.Range("A56").Value = "Korean Statistical Information Services - KOSIS. Statistical Database. Economy / Corporate Business (Company). Census of Establishments. 9th Revision. By province, industrial classification and scale of establishment Period Annual 2006~2012."
$ws.Range("A56").Style = "source"

$ws.Range("A57").Value = ""
$ws.Range("A57").Style = "source"

$ws.Range("A58").Value = "http://kosis.kr/eng/"
$ws.Range("A58").Style = "source"

$ws.Range("A61").Value = "KOSIS"
$ws.Range("A61").Style = "title"

$ws.Range("A62").Value = "Small Business Corportaion available at https://www.sbc.or.kr/sbc/eng/smes/definition.jsp. Article 2 of Framework Act on SMEs and Article 3 of Enforcement Decree of the Act. For micro-enterprises, Article 2 of the Act of Special Measures on Assisting Small Business and Micro-enterprises shall apply."
$ws.Range("A62").Style = "source"
